# Update vm_pu.xlsx results for the "380 kV" case run: voltage magnitudes
# (per-unit) for buses 0-23 (rows 2-25), columns B-F and I-N, change from
# the previous 1.05 pu slack-bus setpoint run to the new 1.02 pu run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataBF = New-Object 'object[,]' 24,5
$dataBF[0,0] = 1.02
$dataBF[0,1] = 1.045818994846036
$dataBF[0,2] = 1.053450030206506
$dataBF[0,3] = 1.049485684645299
$dataBF[0,4] = 1.06078721916428
$dataBF[1,0] = 1.02
$dataBF[1,1] = 1.046832255539727
$dataBF[1,2] = 1.054287329562006
$dataBF[1,3] = 1.050448200525592
$dataBF[1,4] = 1.061804289798447
$dataBF[2,0] = 1.02
$dataBF[2,1] = 1.047488142222786
$dataBF[2,2] = 1.054829292733861
$dataBF[2,3] = 1.051071583725149
$dataBF[2,4] = 1.062463019932425
$dataBF[3,0] = 1.02
$dataBF[3,1] = 1.04776393397886
$dataBF[3,2] = 1.055057175052615
$dataBF[3,3] = 1.051333790048752
$dataBF[3,4] = 1.062740097206799
$dataBF[4,0] = 1.02
$dataBF[4,1] = 1.047810243962268
$dataBF[4,2] = 1.055095439878514
$dataBF[4,3] = 1.051377823600984
$dataBF[4,4] = 1.062786628308013
$dataBF[5,0] = 1.02
$dataBF[5,1] = 1.047491827142123
$dataBF[5,2] = 1.05483233754842
$dataBF[5,3] = 1.051075086805243
$dataBF[5,4] = 1.062466721675276
$dataBF[6,0] = 1.02
$dataBF[6,1] = 1.046161381329088
$dataBF[6,2] = 1.053732962846859
$dataBF[6,3] = 1.049810852737499
$dataBF[6,4] = 1.061130814692008
$dataBF[7,0] = 1.02
$dataBF[7,1] = 1.043818810576075
$dataBF[7,2] = 1.051797094599239
$dataBF[7,3] = 1.047587514577871
$dataBF[7,4] = 1.058781534726635
$dataBF[8,0] = 1.02
$dataBF[8,1] = 1.042258351204235
$dataBF[8,2] = 1.050507481362008
$dataBF[8,3] = 1.046108288764493
$dataBF[8,4] = 1.057218586232545
$dataBF[9,0] = 1.02
$dataBF[9,1] = 1.041582954789049
$dataBF[9,2] = 1.049949302410095
$dataBF[9,3] = 1.045468486697782
$dataBF[9,4] = 1.056542587816229
$dataBF[10,0] = 1.02
$dataBF[10,1] = 1.041332126552093
$dataBF[10,2] = 1.049742005407496
$dataBF[10,3] = 1.045230943303254
$dataBF[10,4] = 1.056291607906448
$dataBF[11,0] = 1.02
$dataBF[11,1] = 1.041385928065148
$dataBF[11,2] = 1.049786469700851
$dataBF[11,3] = 1.04528189229083
$dataBF[11,4] = 1.056345438691282
$dataBF[12,0] = 1.02
$dataBF[12,1] = 1.041562220349398
$dataBF[12,2] = 1.049932166447009
$dataBF[12,3] = 1.045448849088957
$dataBF[12,4] = 1.056521839369147
$dataBF[13,0] = 1.02
$dataBF[13,1] = 1.041670845685305
$dataBF[13,2] = 1.050021939746025
$dataBF[13,3] = 1.04555173095415
$dataBF[13,4] = 1.056630541035818
$dataBF[14,0] = 1.02
$dataBF[14,1] = 1.042303181215015
$dataBF[14,2] = 1.050544530805567
$dataBF[14,3] = 1.046150765377485
$dataBF[14,4] = 1.057263466259537
$dataBF[15,0] = 1.02
$dataBF[15,1] = 1.042699907026806
$dataBF[15,2] = 1.050872400962259
$dataBF[15,3] = 1.046526714974222
$dataBF[15,4] = 1.057660689677762
$dataBF[16,0] = 1.02
$dataBF[16,1] = 1.042931338725492
$dataBF[16,2] = 1.051063664427628
$dataBF[16,3] = 1.046746068720383
$dataBF[16,4] = 1.057892457499925
$dataBF[17,0] = 1.02
$dataBF[17,1] = 1.043010255742585
$dataBF[17,2] = 1.051128884098
$dataBF[17,3] = 1.046820874316719
$dataBF[17,4] = 1.057971496894342
$dataBF[18,0] = 1.02
$dataBF[18,1] = 1.042657339179653
$dataBF[18,2] = 1.050837221280889
$dataBF[18,3] = 1.046486372029796
$dataBF[18,4] = 1.057618063693513
$dataBF[19,0] = 1.02
$dataBF[19,1] = 1.041510305448655
$dataBF[19,2] = 1.04988926140115
$dataBF[19,3] = 1.045399681499115
$dataBF[19,4] = 1.056469890559743
$dataBF[20,0] = 1.02
$dataBF[20,1] = 1.040789374795965
$dataBF[20,2] = 1.049293447187443
$dataBF[20,3] = 1.044717058726441
$dataBF[20,4] = 1.055748659948155
$dataBF[21,0] = 1.02
$dataBF[21,1] = 1.041171529579947
$dataBF[21,2] = 1.049609279853797
$dataBF[21,3] = 1.045078870769029
$dataBF[21,4] = 1.056130934122794
$dataBF[22,0] = 1.02
$dataBF[22,1] = 1.042676573670085
$dataBF[22,2] = 1.050853117395551
$dataBF[22,3] = 1.046504601057524
$dataBF[22,4] = 1.057637324311922
$dataBF[23,0] = 1.02
$dataBF[23,1] = 1.044424201209699
$dataBF[23,2] = 1.052297395987201
$dataBF[23,3] = 1.048161774923931
$dataBF[23,4] = 1.059388311995223

$ws.Range("B2:F25").Value = $dataBF

$dataIN = New-Object 'object[,]' 24,6
$dataIN[0,0] = 1.046568187511246
$dataIN[0,1] = 1.050876694461945
$dataIN[0,2] = 1.056195965933123
$dataIN[0,3] = 1.052242610569663
$dataIN[0,4] = 1.063513048958521
$dataIN[0,5] = 1.052369059905323
$dataIN[1,0] = 1.046890954876427
$dataIN[1,1] = 1.051537457309625
$dataIN[1,2] = 1.056846250181206
$dataIN[1,3] = 1.053016994076385
$dataIN[1,4] = 1.06434410218589
$dataIN[1,5] = 1.053030761112036
$dataIN[2,0] = 1.04709859947846
$dataIN[2,1] = 1.05196464834423
$dataIN[2,2] = 1.057266554477654
$dataIN[2,3] = 1.053518022070909
$dataIN[2,4] = 1.064881854305419
$dataIN[2,5] = 1.05345855880691
$dataIN[3,0] = 1.047185603889778
$dataIN[3,1] = 1.052144151184804
$dataIN[3,2] = 1.057443136605857
$dataIN[3,3] = 1.053728641729771
$dataIN[3,4] = 1.065107925917149
$dataIN[3,5] = 1.053638316562107
$dataIN[4,0] = 1.047200195334978
$dataIN[4,1] = 1.052174285311496
$dataIN[4,2] = 1.05747277883341
$dataIN[4,3] = 1.053764004947439
$dataIN[4,4] = 1.06514588435445
$dataIN[4,5] = 1.053668493482717
$dataIN[5,0] = 1.047099763173109
$dataIN[5,1] = 1.05196704721741
$dataIN[5,2] = 1.057268914424513
$dataIN[5,3] = 1.053520836430953
$dataIN[5,4] = 1.064884875083671
$dataIN[5,5] = 1.053460961086765
$dataIN[6,0] = 1.046677518080756
$dataIN[6,1] = 1.051100078001845
$dataIN[6,2] = 1.056415830173131
$dataIN[6,3] = 1.052504327099709
$dataIN[6,4] = 1.063793905866713
$dataIN[6,5] = 1.052592760675472
$dataIN[7,0] = 1.045924234101305
$dataIN[7,1] = 1.049569579858171
$dataIN[7,2] = 1.054908986951454
$dataIN[7,3] = 1.050712750018922
$dataIN[7,4] = 1.061871544266581
$dataIN[7,5] = 1.051060089048883
$dataIN[8,0] = 1.04541585353284
$dataIN[8,1] = 1.048547400237252
$dataIN[8,2] = 1.05390203478608
$dataIN[8,3] = 1.049518156159177
$dataIN[8,4] = 1.060590047773212
$dataIN[8,5] = 1.050036457815656
$dataIN[9,0] = 1.045194253953326
$dataIN[9,1] = 1.048104352088887
$dataIN[9,2] = 1.053465452383665
$dataIN[9,3] = 1.049000840295394
$dataIN[9,4] = 1.060035170675554
$dataIN[9,5] = 1.049592780488103
$dataIN[10,0] = 1.045111721745315
$dataIN[10,1] = 1.04793971887524
$dataIN[10,2] = 1.053303201690852
$dataIN[10,3] = 1.048808679231418
$dataIN[10,4] = 1.059829067975644
$dataIN[10,5] = 1.049427913476409
$dataIN[11,0] = 1.045129435147836
$dataIN[11,1] = 1.047975036212583
$dataIN[11,2] = 1.053338008834754
$dataIN[11,3] = 1.048849898738287
$dataIN[11,4] = 1.05987327754504
$dataIN[11,5] = 1.049463280968422
$dataIN[12,0] = 1.045187436307235
$dataIN[12,1] = 1.048090744788327
$dataIN[12,2] = 1.053452042413218
$dataIN[12,3] = 1.04898495634061
$dataIN[12,4] = 1.060018134090905
$dataIN[12,5] = 1.049579153863615
$dataIN[13,0] = 1.045223143555855
$dataIN[13,1] = 1.048162028004286
$dataIN[13,2] = 1.053522291073004
$dataIN[13,3] = 1.049068168882641
$dataIN[13,4] = 1.060107385448555
$dataIN[13,5] = 1.049650538309917
$dataIN[14,0] = 1.045430529467616
$dataIN[14,1] = 1.048576794692722
$dataIN[14,2] = 1.05393099744307
$dataIN[14,3] = 1.049552487716037
$dataIN[14,4] = 1.060626873584957
$dataIN[14,5] = 1.050065894014625
$dataIN[15,0] = 1.045560224442585
$dataIN[15,1] = 1.048836849955412
$dataIN[15,2] = 1.054187217044756
$dataIN[15,3] = 1.049856275344191
$dataIN[15,4] = 1.060952740266439
$dataIN[15,5] = 1.050326318585617
$dataIN[16,0] = 1.045635731726555
$dataIN[16,1] = 1.04898849354251
$dataIN[16,2] = 1.054336611016579
$dataIN[16,3] = 1.050033464787528
$dataIN[16,4] = 1.061142814624026
$dataIN[16,5] = 1.050478177524002
$dataIN[17,0] = 1.045661453734867
$dataIN[17,1] = 1.049040192910489
$dataIN[17,2] = 1.054387541229073
$dataIN[17,3] = 1.050093880991936
$dataIN[17,4] = 1.061207625365541
$dataIN[17,5] = 1.050529950311013
$dataIN[18,0] = 1.045546324044068
$dataIN[18,1] = 1.048808952854212
$dataIN[18,2] = 1.054159732749528
$dataIN[18,3] = 1.049823682293948
$dataIN[18,4] = 1.060917777662283
$dataIN[18,5] = 1.050298381867333
$dataIN[19,0] = 1.045170362489929
$dataIN[19,1] = 1.048056673299446
$dataIN[19,2] = 1.053418464697291
$dataIN[19,3] = 1.048945185430563
$dataIN[19,4] = 1.059975477352189
$dataIN[19,5] = 1.049545033989312
$dataIN[20,0] = 1.044932706302372
$dataIN[20,1] = 1.047583307046346
$dataIN[20,2] = 1.052951911075623
$dataIN[20,3] = 1.048392799764735
$dataIN[20,4] = 1.059383035414353
$dataIN[20,5] = 1.049070995501836
$dataIN[21,0] = 1.045058813056612
$dataIN[21,1] = 1.047834283190828
$dataIN[21,2] = 1.053199286191327
$dataIN[21,3] = 1.048685633491785
$dataIN[21,4] = 1.059697098078495
$dataIN[21,5] = 1.049322328061231
$dataIN[22,0] = 1.045552605473313
$dataIN[22,1] = 1.048821558483928
$dataIN[22,2] = 1.054172151888129
$dataIN[22,3] = 1.049838409705014
$dataIN[22,4] = 1.060933575752447
$dataIN[22,5] = 1.050311005398489
$dataIN[23,0] = 1.046120068589367
$dataIN[23,1] = 1.049965577914344
$dataIN[23,2] = 1.055298965520645
$dataIN[23,3] = 1.051175954804926
$dataIN[23,4] = 1.062368509570562
$dataIN[23,5] = 1.051456649467717

$ws.Range("I2:N25").Value = $dataIN

Write-Host "Done"